$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights_and_functions")

$ws.Cells.Item(100,2).Value = "hello"
$ws.Cells.Item(61,2).Copy()
$ws.Cells.Item(100,2).PasteSpecial(-4122)  # xlPasteFormats = -4122
Write-Host "done"
